# Add a new header row (1:1) above the existing data, pushing everything
# down by one row, then populate the two new header cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("1:1").Insert()

$ws.Range("A1").Value = "Female"
$ws.Range("B1").Value = "Chainese"

# Restore the sort range / sort state that shifted along with the data
# (F:I used to be sorted F2:I41 on F2:F41; after the insert the same data
# now lives one row down, at F3:I42 / F3:F42).
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("F3:F42"))
$sortObj.SetRange($ws.Range("F3:I42"))
$sortObj.Header = -4142
$sortObj.Apply()

# Match the saved selection/active cell from the edited workbook.
$ws.Range("B2").Select()
